# Commit: Change Xianfeng Zhang to Xianfeng Zeng (#288)
# The "Researcher Name" column (C) on the Samples sheet lists the same
# researcher, "Xianfeng Zhang", for every sample row (C2:C17). Correct the
# spelling to "Xianfeng Zeng" everywhere it appears, in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

$ws.Range("C2:C17").Replace("Xianfeng Zhang", "Xianfeng Zeng")
